$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3:F9").Value = 0
